# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 225
$ws.Range("I6").Value = 142.5
$ws.Range("K6").Value = 427.5
$ws.Range("M6").Value = -315.5

$ws.Range("H12").Value = 192.14285
$ws.Range("I12").Value = 229.2
$ws.Range("J12").Value = 99.5
$ws.Range("K12").Value = 229.2
$ws.Range("L12").Value = 99.5
$ws.Range("M12").Value = -59.19999999999999
$ws.Range("N12").Value = -439.5

$ws.Range("H98").Value = 3174.7222
$ws.Range("I98").Value = 2973.8462
$ws.Range("J98").Value = 3697
$ws.Range("K98").Value = 2973.8462
$ws.Range("L98").Value = 3697
$ws.Range("M98").Value = -1475.8462
$ws.Range("N98").Value = -6693

$ws.Range("H106").Value = 4948.75
$ws.Range("I106").Value = 4948.75
$ws.Range("K106").Value = 4948.75
$ws.Range("M106").Value = -4317.75

$ws.Range("H122").Value = 3174.7222
$ws.Range("I122").Value = 2973.8462
$ws.Range("J122").Value = 3697
$ws.Range("K122").Value = 8921.5386
$ws.Range("L122").Value = 11091
$ws.Range("M122").Value = -6471.5386
$ws.Range("N122").Value = -15991

$ws.Range("H131").Value = 1382
$ws.Range("J131").Value = 3300
$ws.Range("L131").Value = 9900
$ws.Range("N131").Value = -19980

$ws.Range("H137").Value = 1322.2174
$ws.Range("I137").Value = 1151.25
$ws.Range("J137").Value = 2462
$ws.Range("K137").Value = 3453.75
$ws.Range("L137").Value = 7386
$ws.Range("M137").Value = -903.75
$ws.Range("N137").Value = -12486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1895
$ws.Range("I45").Value = 1895
$ws.Range("K45").Value = 1895
$ws.Range("M45").Value = -1518

$ws.Range("H132").Value = 1879.8235
$ws.Range("I132").Value = 1863.8667
$ws.Range("K132").Value = 5591.6001
$ws.Range("M132").Value = -3061.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1333.8889
$ws.Range("I105").Value = 1333.8889
$ws.Range("K105").Value = 1333.8889
$ws.Range("M105").Value = 413.1111000000001

$ws.Range("H107").Value = 1297.5714
$ws.Range("I107").Value = 680.5
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 680.5
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 1239.5
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2704
$ws.Range("I31").Value = 3143.3333
$ws.Range("J31").Value = 2411.111
$ws.Range("K31").Value = 3143.3333
$ws.Range("L31").Value = 2411.111
$ws.Range("M31").Value = -2848.3333
$ws.Range("N31").Value = -3001.111

$ws.Range("H34").Value = 2704
$ws.Range("I34").Value = 3143.3333
$ws.Range("J34").Value = 2411.111
$ws.Range("K34").Value = 3143.3333
$ws.Range("L34").Value = 2411.111
$ws.Range("M34").Value = -2941.3333
$ws.Range("N34").Value = -2815.111

$ws.Range("H58").Value = 2543.7917
$ws.Range("I58").Value = 2321.5789
$ws.Range("K58").Value = 2321.5789
$ws.Range("M58").Value = -2118.5789

$ws.Range("H94").Value = 2164
$ws.Range("I94").Value = 2299.6
$ws.Range("K94").Value = 2299.6
$ws.Range("M94").Value = -1848.6

$ws.Range("H107").Value = 1120
$ws.Range("I107").Value = 1476.75
$ws.Range("J107").Value = 644.3333
$ws.Range("K107").Value = 1476.75
$ws.Range("L107").Value = 644.3333
$ws.Range("M107").Value = 443.25
$ws.Range("N107").Value = -4484.3333

$ws.Range("H122").Value = 1329.3
$ws.Range("I122").Value = 1141.8572
$ws.Range("J122").Value = 1766.6666
$ws.Range("K122").Value = 3425.5716
$ws.Range("L122").Value = 5299.9998
$ws.Range("M122").Value = -975.5715999999998
$ws.Range("N122").Value = -10199.9998

$ws.Range("H136").Value = 2543.7917
$ws.Range("I136").Value = 2321.5789
$ws.Range("K136").Value = 6964.736699999999
$ws.Range("M136").Value = -4414.736699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 455
$ws.Range("J12").Value = 459.6
$ws.Range("L12").Value = 1378.8
$ws.Range("N12").Value = -1724.8

$ws.Range("H121").Value = 787.8889
$ws.Range("J121").Value = 999.4286
$ws.Range("L121").Value = 2998.2858
$ws.Range("N121").Value = -5618.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 508.36365
$ws.Range("I2").Value = 591.4286
$ws.Range("K2").Value = 591.4286
$ws.Range("M2").Value = -478.4286

$ws.Range("H70").Value = 3333
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3333
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3333
$ws.Range("N70").Value = -3873
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 3333
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3333
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3333
$ws.Range("N73").Value = -5205
$ws.Range("M73").ClearContents()

$ws.Range("H80").Value = 3995.2173
$ws.Range("I80").Value = 2635
$ws.Range("J80").Value = 4199.25
$ws.Range("K80").Value = 2635
$ws.Range("L80").Value = 4199.25
$ws.Range("M80").Value = -1637
$ws.Range("N80").Value = -6195.25

$ws.Range("H83").Value = 3995.2173
$ws.Range("I83").Value = 2635
$ws.Range("J83").Value = 4199.25
$ws.Range("K83").Value = 13175
$ws.Range("L83").Value = 20996.25
$ws.Range("M83").Value = -8183
$ws.Range("N83").Value = -30980.25

$ws.Range("H99").Value = 27177.889
$ws.Range("I99").Value = 27177.889
$ws.Range("K99").Value = 27177.889
$ws.Range("M99").Value = -24931.889

$ws.Range("H102").Value = 546.05
$ws.Range("I102").Value = 522.1579
$ws.Range("K102").Value = 522.1579
$ws.Range("M102").Value = 1099.8421

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1568.5264
$ws.Range("I16").Value = 1738.2307
$ws.Range("K16").Value = 1738.2307
$ws.Range("M16").Value = -1568.2307

$ws.Range("H40").Value = 3167
$ws.Range("I40").Value = 2656.7693
$ws.Range("K40").Value = 2656.7693
$ws.Range("M40").Value = -2520.7693

$ws.Range("H55").Value = 1800
$ws.Range("I55").Value = 1800
$ws.Range("K55").Value = 1800
$ws.Range("M55").Value = -1627

$ws.Range("H107").Value = 3500
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580

$ws.Range("H122").Value = 4065
$ws.Range("I122").Value = 3502.8572
$ws.Range("K122").Value = 10508.5716
$ws.Range("M122").Value = -8058.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16107.2
$ws.Range("I41").Value = 16897.75
$ws.Range("J41").Value = 15580.167
$ws.Range("K41").Value = 16897.75
$ws.Range("L41").Value = 15580.167
$ws.Range("M41").Value = -16507.75
$ws.Range("N41").Value = -16360.167

$ws.Range("H113").Value = 1216.0714
$ws.Range("I113").Value = 1125.0769
$ws.Range("K113").Value = 3375.2307
$ws.Range("M113").Value = -1205.2307

$ws.Range("H122").Value = 1222
$ws.Range("I122").Value = 1222
$ws.Range("K122").Value = 3666
$ws.Range("M122").Value = -1216
